# New sheet made for December
# The existing sheet is being cleared/repurposed for the new month: some
# "x"/"X" marks are removed and a new "date gone" is recorded for Madeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Eric (row 3): clear "Closing (tear down)" and "Facebook Events" marks
$ws.Range("G3").ClearContents()
$ws.Range("I3").ClearContents()

# Madeline (row 10): record a new date she is gone - Dec 5, 2019
$ws.Range("M10").Value = "12/5/2019"

# Jessica (row 11): clear "Opening", "First Door Shift", "Teaching (follow)"
# and "Facebook Events" marks
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("I11").ClearContents()

# Move the active selection to L13
$ws.Range("L13").Select()
